# Apply the crypto price/volume refresh for Mon May  8 14:46:44 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Ref = 'D2'; Value = '27.918.60' },
    @{ Ref = 'D3'; Value = '1.863.84' },
    @{ Ref = 'E3'; Value = '  -3.23%  ' },
    @{ Ref = 'D4'; Value = '1.004' },
    @{ Ref = 'E4'; Value = '  +0.33%  ' },
    @{ Ref = 'D5'; Value = '318.46' },
    @{ Ref = 'E6'; Value = '  +0.38%  ' },
    @{ Ref = 'D7'; Value = '0.4365' },
    @{ Ref = 'E7'; Value = '  -5.16%  ' },
    @{ Ref = 'D8'; Value = '0.3732' },
    @{ Ref = 'E8'; Value = '  -2.61%  ' },
    @{ Ref = 'D9'; Value = '0.07493' },
    @{ Ref = 'E9'; Value = '  -3.26%  ' },
    @{ Ref = 'D10'; Value = '0.9364' },
    @{ Ref = 'E10'; Value = '  -4.82%  ' },
    @{ Ref = 'D11'; Value = '21.27' },
    @{ Ref = 'E11'; Value = '  -5.01%  ' },
    @{ Ref = 'D12'; Value = '1.884.65' },
    @{ Ref = 'E12'; Value = '  -2.55%  ' },
    @{ Ref = 'D13'; Value = '6.743' },
    @{ Ref = 'E13'; Value = '  -3.38%  ' },
    @{ Ref = 'D14'; Value = '5.446' },
    @{ Ref = 'E14'; Value = '  -4.48%  ' },
    @{ Ref = 'D15'; Value = '0.06899' },
    @{ Ref = 'E15'; Value = '  -1.97%  ' },
    @{ Ref = 'D16'; Value = '1.006' },
    @{ Ref = 'D17'; Value = '81.45' },
    @{ Ref = 'E17'; Value = '  -3.45%  ' },
    @{ Ref = 'D18'; Value = '0.000009055' },
    @{ Ref = 'E18'; Value = '  -4.99%  ' },
    @{ Ref = 'D19'; Value = '1.005' },
    @{ Ref = 'E19'; Value = '  +0.39%  ' },
    @{ Ref = 'D20'; Value = '15.89' },
    @{ Ref = 'E20'; Value = '  -5.05%  ' },
    @{ Ref = 'D21'; Value = '27.896.83' },
    @{ Ref = 'E21'; Value = '  -3.91%  ' },
    @{ Ref = 'D22'; Value = '5.127' },
    @{ Ref = 'E22'; Value = '  -4.10%  ' },
    @{ Ref = 'D23'; Value = '11.06' },
    @{ Ref = 'E23'; Value = '  +0.72%  ' },
    @{ Ref = 'D24'; Value = '2.159.44' },
    @{ Ref = 'E24'; Value = '  +0.09%  ' },
    @{ Ref = 'D25'; Value = '2.011' },
    @{ Ref = 'E25'; Value = '  -3.96%  ' },
    @{ Ref = 'D26'; Value = '154.71' },
    @{ Ref = 'E26'; Value = '  -2.21%  ' },
    @{ Ref = 'D27'; Value = '18.50' },
    @{ Ref = 'E27'; Value = '  -3.27%  ' },
    @{ Ref = 'D28'; Value = '5.550' },
    @{ Ref = 'E28'; Value = '  -2.65%  ' },
    @{ Ref = 'D29'; Value = '113.61' },
    @{ Ref = 'E29'; Value = '  -3.73%  ' },
    @{ Ref = 'D30'; Value = '1.712' },
    @{ Ref = 'E30'; Value = '  -7.91%  ' },
    @{ Ref = 'D31'; Value = '0.09032' },
    @{ Ref = 'E31'; Value = '  -3.35%  ' },
    @{ Ref = 'D32'; Value = '0.8164' },
    @{ Ref = 'E32'; Value = '  -5.91%  ' },
    @{ Ref = 'D33'; Value = '4.826' },
    @{ Ref = 'E33'; Value = '  -5.82%  ' },
    @{ Ref = 'E34'; Value = '  -6.61%  ' },
    @{ Ref = 'D35'; Value = '2.971' },
    @{ Ref = 'E35'; Value = '  -2.55%  ' },
    @{ Ref = 'D36'; Value = '1.005' },
    @{ Ref = 'E36'; Value = '  +0.40%  ' },
    @{ Ref = 'D37'; Value = '0.05527' },
    @{ Ref = 'E37'; Value = '  -3.27%  ' },
    @{ Ref = 'E38'; Value = '  -3.07%  ' },
    @{ Ref = 'E39'; Value = '  -3.65%  ' },
    @{ Ref = 'D40'; Value = '2.948' },
    @{ Ref = 'E40'; Value = '  -3.29%  ' },
    @{ Ref = 'D41'; Value = '0.5270' },
    @{ Ref = 'E42'; Value = '  -6.82%  ' },
    @{ Ref = 'D43'; Value = '0.1701' },
    @{ Ref = 'E43'; Value = '  -3.08%  ' },
    @{ Ref = 'D44'; Value = '8.802' },
    @{ Ref = 'E44'; Value = '  -6.47%  ' },
    @{ Ref = 'D45'; Value = '0.06751' },
    @{ Ref = 'E45'; Value = '  -2.24%  ' },
    @{ Ref = 'D46'; Value = '0.4899' },
    @{ Ref = 'E46'; Value = '  -5.92%  ' },
    @{ Ref = 'D47'; Value = '10.59' },
    @{ Ref = 'E47'; Value = '  -5.80%  ' },
    @{ Ref = 'D48'; Value = '107.85' },
    @{ Ref = 'B49'; Value = 'RenderToken' },
    @{ Ref = 'C49'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Ref = 'D49'; Value = '1.916' },
    @{ Ref = 'E49'; Value = '  -13.90%  ' },
    @{ Ref = 'B50'; Value = 'PaxDollar' },
    @{ Ref = 'C50'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Ref = 'D50'; Value = '1.004' },
    @{ Ref = 'E50'; Value = '  +0.30%  ' },
    @{ Ref = 'B51'; Value = 'NEARProtocol' },
    @{ Ref = 'C51'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Ref = 'D51'; Value = '1.677' },
    @{ Ref = 'E51'; Value = '  -5.96%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Ref)
    if ($u.Ref.StartsWith("D")) {
        # Column D values look numeric (e.g. "5.550", "0.000009055") but the source
        # data must stay as literal text so trailing zeros / exponent-prone values
        # are preserved exactly as scraped. Force text format, assign, then drop the
        # number-format override so no extra style is left attached to the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
